# Apply the cryptos.xlsx data refresh described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "59.822.16"
$ws.Range("E2").Value = "  +0.15%  "
$ws.Range("D3").Value = "2.532.19"
$ws.Range("E3").Value = "  +1.49%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "544.38"
$ws.Range("E5").Value = "  +0.08%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.95"
$ws.Range("E6").Value = "  -0.87%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.996"
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("E8").Value = "  -1.15%  "
$ws.Range("D9").Value = "2.564.39"
$ws.Range("E9").Value = "  +1.58%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.101"
$ws.Range("E10").Value = "  +0.39%  "
$ws.Range("E11").Value = "  +0.64%  "
$ws.Range("E12").Value = "  +3.21%  "
$ws.Range("E13").Value = "  +1.06%  "
$ws.Range("D14").Value = "2.980.92"
$ws.Range("E14").Value = "  +1.39%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "23.70"
$ws.Range("E15").Value = "  -4.01%  "
$ws.Range("D16").Value = "59.786.20"
$ws.Range("E16").Value = "  -0.27%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000143"
$ws.Range("E17").Value = "  +1.70%  "
$ws.Range("D18").Value = "2.545.25"
$ws.Range("E18").Value = "  +1.38%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.27"
$ws.Range("E19").Value = "  -2.15%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.32"
$ws.Range("E20").Value = "  -1.43%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "328.10"
$ws.Range("E21").Value = "  +0.09%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.998"
$ws.Range("E22").Value = "  +0.12%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.96"
$ws.Range("E23").Value = "  +2.51%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "62.35"
$ws.Range("E24").Value = "  +1.25%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.440"
$ws.Range("E25").Value = "  -2.13%  "
$ws.Range("E26").Value = "  +2.07%  "
$ws.Range("E27").Value = "  -1.40%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.05"
$ws.Range("E28").Value = "  +2.07%  "
$ws.Range("B29").Value = "Aptos"
$ws.Range("C29").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.98"
$ws.Range("E29").Value = "  +0.78%  "
$ws.Range("B30").Value = "PEPE"
$ws.Range("C30").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D30").Value = "0.0₃0800"
$ws.Range("E30").Value = "  +0.55%  "
$ws.Range("E31").Value = "  -0.02%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.23"
$ws.Range("E32").Value = "  -6.06%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "161.61"
$ws.Range("E33").Value = "  +1.77%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.47"
$ws.Range("E34").Value = "  +1.67%  "
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.80"
$ws.Range("E36").Value = "  -0.74%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.46"
$ws.Range("E37").Value = "  -1.77%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.63"
$ws.Range("E38").Value = "  -6.71%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.70"
$ws.Range("E39").Value = "  -6.53%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "37.16"
$ws.Range("E40").Value = "  +1.11%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "304.03"
$ws.Range("E41").Value = "  -4.12%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.845"
$ws.Range("E42").Value = "  +1.40%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.73"
$ws.Range("E43").Value = "  -1.48%  "
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.994"
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("B45").Value = "Mantle"
$ws.Range("C45").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.609"
$ws.Range("E45").Value = "  +0.33%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.81"
$ws.Range("E46").Value = "  +0.24%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "19.10"
$ws.Range("E47").Value = "  +2.12%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0938"
$ws.Range("E48").Value = "  -0.72%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "123.56"
$ws.Range("E49").Value = "  -2.97%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0519"
$ws.Range("E50").Value = "  -2.75%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0229"
$ws.Range("E51").Value = "  -1.33%  "
